# Fruta / hortaliza, semanal
# A new weekly observation is inserted as row 21 ("Mora" at Mercado Mayorista
# Lo Valledor de Santiago, dated 2021-12-10 / serial 44540), which pushes the
# previously existing rows 21-47 down to rows 22-48 (Excel carries the D-column
# date style along automatically when inserting a whole row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 21; everything below (old rows 21-47)
# shifts down to 22-48, and the sheet's used range grows to A1:T48.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new observation.
$ws.Range("A21").Value = 6
$ws.Range("B21").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value2 = 44540
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100101
$ws.Range("H21").Value = "Berries"
$ws.Range("I21").Value = 100101008
$ws.Range("J21").Value = "Mora"
$ws.Range("K21").Value = "Sin especificar"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 250
$ws.Range("N21").Value = 5000
$ws.Range("O21").Value = 5000
$ws.Range("P21").Value = 5000
$ws.Range("Q21").Value = "$/bandeja 2 kilos"
$ws.Range("R21").Value = "Provincia de Curicó"
$ws.Range("S21").Value = 2500
$ws.Range("T21").Value = 2
